$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Tên món"
$ws.Range("B1").Value = "Đơn vị tính"
$ws.Range("C1").Value = "Giá gốc"
$ws.Range("D1").Value = "Giá bán"

# Data rows
$data = @(
    @("Cà phê sữa", "Ly", 8000, 9000),
    @("Cà phê đen", "Ly", 7000, 8000),
    @("Thuốc White Horse", "Gói", 20000, 24000),
    @("Nước cam", "Ly", 5000, 8000),
    @("Chanh dây", "Ly", 12000, 18000),
    @("Coca Cola", "Chai", 3000, 7000),
    @("Number One", "Chai", 7000, 10000),
    @("Bò húc", "Lon", 12000, 18000),
    @("Nutriboost", "Chai", 12000, 18000),
    @("Nước dừa", "Trái", 30000, 50000),
    @("Lemoncello", "Shot", 70000, 90000)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}

$ws.Range("A13").Select()
